$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 42731.0647474537

$ws.Range("B3").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("D3").Value = $newDate
$ws.Range("E3").Value = $newDate
$ws.Range("F3").Value = $newDate
$ws.Range("G3").Value = $newDate
